# Update countries & provincias Spain
# Applies the data refresh for the "Pais" worksheet:
#   - Update the "Datos actualizados..." timestamp string (row 1)
#   - Update a handful of per-country statistic rows with new figures

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: refreshed timestamp text
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 07:17"

# Row 5 (India)
$ws.Range("B5").Value = 4204613
$ws.Range("C5").Value = 2051
$ws.Range("D5").Value = 3250429
$ws.Range("E5").Value = 882497

# Row 62 (Kirguistan)
$ws.Range("B62").Value = 44458
$ws.Range("C62").Value = 55
$ws.Range("D62").Value = 39960
$ws.Range("E62").Value = 3438

# Row 64 (Uzbekistan)
$ws.Range("B64").Value = 43663
$ws.Range("C64").Value = 76
$ws.Range("D64").Value = 41277
$ws.Range("E64").Value = 2038
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 348

# Row 125 (Tailandia)
$ws.Range("B125").Value = 3445
$ws.Range("C125").Value = 1
$ws.Range("E125").Value = 106

# Row 187 (Butan)
$ws.Range("B187").Value = 230
$ws.Range("C187").Value = 2
$ws.Range("D187").Value = 151
$ws.Range("E187").Value = 79
